$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEDG")

# Row 4 - Inventory
$ws.Range("B4").Value = 332000000.0
$ws.Range("C4").Value = 297000000.0
$ws.Range("D4").Value = 264000000.0
$ws.Range("E4").Value = 199000000.0
$ws.Range("F4").Value = 171000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 162000000.0
$ws.Range("C15").Value = 122000000.0
$ws.Range("D15").Value = 158000000.0
$ws.Range("E15").Value = 139000000.0
$ws.Range("F15").Value = 157000000.0

# Row 25 - Long Term Tax Liability (Deferred)
$ws.Range("B25").Value = -3000000.0
$ws.Range("C25").Value = -5000000.0
$ws.Range("D25").Value = -18000000.0
$ws.Range("E25").Value = -15000000.0
$ws.Range("F25").Value = -12000000.0
